$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" hidden bookmark (currently sitting at
#    the start of the "Υλικός Εξοπλισμός" TOC entry). Word will
#    renumber the remaining bookmarks (the "_Hlk526270059" bookmark
#    shifts from id 3 down to id 2) automatically on save.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Re-create "_GoBack" at the new location: right after the page
#    break run that separates "ΚΑΤΑΛΟΓΟΣ ΕΙΚΟΝΩΝ" from
#    "ΚΑΤΑΛΟΓΟΣ ΠΙΝΑΚΩΝ" (i.e. Word moved the "last edit" marker here
#    once chapter 5 was finished).
#
#    Find the (now empty) paragraph that contains only the page
#    break, right before the "ΚΑΤΑΛΟΓΟΣ ΠΙΝΑΚΩΝ" heading.
# ------------------------------------------------------------------
$target = $null
$paras = $d.Paragraphs
$pageBreakOnly = [string]([char]12) + [string]([char]13)
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq $pageBreakOnly) {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.StartsWith("ΚΑΤΑΛΟΓΟΣ ΠΙΝΑΚΩΝ")) {
            $target = $p
        }
    }
}

if ($target -ne $null) {
    $insertPos = $target.Range.End - 1

    # The COM layer mis-resolves a *collapsed* range sitting exactly on
    # the paragraph's content/mark boundary, so nudge it into a safe
    # spot: drop a placeholder character right after the break, anchor
    # the bookmark just in front of it, then remove the placeholder
    # again. The bookmark stays put, now flush against the run.
    $guard = $d.Range($insertPos, $insertPos)
    $guard.InsertAfter("X")

    $bmRange = $d.Range($insertPos, $insertPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $d.Range($insertPos, $insertPos + 1).Delete()
}
